# "changed report creation a bit"
#
#  - Column B (item code) loses its trailing " xxx-n" suffix, keeping only
#    the leading code (e.g. "EK008104 ekt-4" -> "EK008104").
#  - Column D is reset to 0 for every row.
#  - Column E becomes the literal text "0%" for every row (was a numeric
#    percentage/quantity).
#  - Column F becomes the literal text "specs" for every row (was the
#    placeholder Russian sentence).
# Columns A and C are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, shortened item codes for column B, keyed by row number.
$codes = @{
    1 = "EK008104"
    2 = "EK008104"
    3 = "EK008104"
    4 = "EK008105"
    5 = "EK008105"
    6 = "EK008105"
    7 = "EK008105"
    8 = "MS008101"
    9 = "MS008101"
}

# Force column E to "Text" number format *before* writing "0%" into it so
# Excel stores the literal string instead of re-interpreting it as the
# number 0 formatted as a percentage. Restore the cells' style afterwards
# so nothing looks different from a plain text cell.
$eRange = $ws.Range("E1:E9")
$eRange.NumberFormat = "@"

for ($row = 1; $row -le 9; $row++) {
    $ws.Cells.Item($row, 2).Value = $codes[$row]   # B: shortened code
    $ws.Cells.Item($row, 4).Value = 0               # D: always 0
    $ws.Cells.Item($row, 5).Value = "0%"            # E: literal text "0%"
    $ws.Cells.Item($row, 6).Value = "specs"         # F: literal text "specs"
}

$eRange.Style = "Normal"
